# Regenerate orders with updated distance/size codes.
# D51 -> D55, D64 -> D69, D80 -> D86, S30 -> S31
# (applies everywhere these tokens occur: Condition, Filename_Left,
#  Filename_Right, Distance and Size columns/tables)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$range = $ws.UsedRange

$range.Replace("D51", "D55") | Out-Null
$range.Replace("D64", "D69") | Out-Null
$range.Replace("D80", "D86") | Out-Null
$range.Replace("S30", "S31") | Out-Null
